# Rebuild the "Sanju Samson" innings table: new columns (ownTeam/oppTeam)
# were inserted after "result" and several new match rows were added.
# Final shape: header row + 7 data rows, columns A..K (11 cols).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")

$rows = @(
    ,@(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Rajasthan Royals", "Mumbai Indians", "Sanju Samson †", "54", "31", "4", "3", "174.19")
    ,@(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", "Sanju Samson †", "9", "6", "0", "1", "150.00")
    ,@(" Abu Dhabi", " October 30 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", "Sanju Samson †", "48", "25", "4", "3", "192.00")
    ,@(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", "Sanju Samson †", "36", "26", "3", "1", "138.46")
    ,@(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Sanju Samson †", "1", "4", "0", "0", "25.00")
    ,@(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Rajasthan Royals", "Chennai Super Kings", "Sanju Samson †", "74", "32", "1", "9", "231.25")
    ,@(" Abu Dhabi", " October 19 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Rajasthan Royals", "Chennai Super Kings", "Sanju Samson †", "0", "3", "0", "0", "0.00")
)

# Columns G..K (totalRuns, totalBalls, total4s, total6s, sr) hold digit-only
# text in the source data and must stay text instead of being coerced to
# numbers, so force a text number-format on that block before writing it.
$numericTextCols = @(7, 8, 9, 10, 11)

$lastCol = $headers.Count
$lastRow = $rows.Count + 1

# Clear out anything left over from the previous (smaller) table footprint.
$ws.Cells.Clear()

for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($excelRow, $c)
        if ($numericTextCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}
